$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Target 1
$found1 = $find.Execute("1. Características do Comportamento Empreendedor: Busca de oportunidades e iniciativa. Correr riscos calculados. Exigência de qualidade e eficiência. Persistência. Comprometimento. Busca de informações. Estabelecimento de metas. Monitoramento e planejamento sistemático. Persuasão e rede contatos. Independência e autoconfiança.2. Estratégia, Inovação e Marketing.3. Design Thinking.4. Modelo de Negócios (Business Model Canvas e Lean Startup - Lean Canvas): Problema. Segmento de Clientes. Proposta de Valor Única. Solução. Métricas-Chave. Canais. Estrutura de Custos. Fluxos de Receita. Vantagem Injusta.5. Produto mínimo viável: Ciclo Construir-Mensurar-Aprender. Valor da vida útil do cliente. Prototipação rápida.6. Gestão de processos e Gerenciamento ágil de projetos.7. Plano de Negócios: Marketing, Finanças, Recursos Humanos, Desenvolvimento de Produtos e Tecnologia da Informação e Comunicação.8. Proposta da criação de uma startup, do modelo de negócios ao plano de negócios, incluindo a montagem do produto mínimo viável e uma rodada de PITCH.9. Desenvolvimento de atividade prática extensionista (produção de conteúdo digital sobre empreendedorismo e inovação)10. Visita (viagem didática complementar) a um ambiente de inovação e empreendedorismo (ex. incubadora/aceleradora ou parque tecnológico), para compreender o desenvolvimento dos processos de empreendedorismo e inovação.", $false, $false, $false, $false, $false, $true, 1, $false, "1. Características do Comportamento Empreendedor: Busca de oportunidades e iniciativa. Correr riscos calculados. Exigência de qualidade e eficiência. Persistência. Comprometimento. Busca de informações. Estabelecimento de metas. Monitoramento e planejamento sistemático. Persuasão e rede contatos. Independência e autoconfiança.^l2. Estratégia, Inovação e Marketing.^l3. Design Thinking.^l4. Modelo de Negócios (Business Model Canvas e Lean Startup - Lean Canvas): Problema. Segmento de Clientes. Proposta de Valor Única. Solução. Métricas-Chave. Canais. Estrutura de Custos. Fluxos de Receita. Vantagem Injusta.^l5. Produto mínimo viável: Ciclo Construir-Mensurar-Aprender. Valor da vida útil do cliente. Prototipação rápida.^l6. Gestão de processos e Gerenciamento ágil de projetos.^l7. Plano de Negócios: Marketing, Finanças, Recursos Humanos, Desenvolvimento de Produtos e Tecnologia da Informação e Comunicação.^l8. Proposta da criação de uma startup, do modelo de negócios ao plano de negócios, incluindo a montagem do produto mínimo viável e uma rodada de PITCH.^l9. Desenvolvimento de atividade prática extensionista (produção de conteúdo digital sobre empreendedorismo e inovação)^l10. Visita (viagem didática complementar) a um ambiente de inovação e empreendedorismo (ex. incubadora/aceleradora ou parque tecnológico), para compreender o desenvolvimento dos processos de empreendedorismo e inovação.", 2)
Write-Host "Target 1 found:" $found1

# Target 2
$found2 = $find.Execute("BLANK, Steve Gary. Do Sonho a realização em 4 passos: Estratégias para a criação de empresas de sucesso. Editora Evora. 3ª edição, 2008BLANK, Steve; DORF, Bob. STARTUP: Manual do Empreendedorismo. O guia passo a passo para construir uma grande empresa. Alta Books Editora.  1ª edição, 2014.CECCONELO, Antonio; AJZENTAL, Alberto. A construção do plano de negócios. Ed. Saraiva, 1ª edição, 2008.CHIAVENATO, Idalberto. Empreendedorismo – dando asas ao espírito empreendedor. Ed. Saraiva, 3ª edição, 2008.DOLABELA, Fernando. O Segredo de Luísa. Rio de Janeiro: Sextante, 2008. DORNELAS, Jose. Empreendedorismo: transformando ideias em negócios. Editora Campus. 1ª edição, 2001DORNELAS, Jose. Empreendedorismo na prática. LTC. 3ª edição, 2015DORNELAS, Jose Carlos Assis. Empreendedorismo na prática – mitos e verdades do empreendedor de sucesso. Elsevier/Campus: Rio de Janeiro, 2007. FILION, L. J.; Visão e Relações: Elementos para um Metamodelo da Atividade Empreendedora. International Small Business Journal, 1991. Tradução de Costa, S.R. FILION, L. J.; - O planejamento do seu Sistema de Aprendizagem Empresarial: Identifique uma Visão e Avalie o seu Sistema de Relações. Revista de Administração de Empresas, FGV, São Paulo, jul/set. 1991, pag. 31(3): 63:71. HASHIMOTO, Marcos. Espírito empreendedor nas organizações – aumentando a competitividade através do intraempreendedorismo. São Paulo: Saraiva, 2006. HISRICH, Robert; PETERS, Michael.  Empreendedorismo. 5.ed. - Porto Alegre: Bookman, 2004. OSTERWALDER, Alexander. Inovação Em Modelos de Negócios – Business Model Generation. Editora Alta Books, 2011PINCHOT, Gifford; PELLMAN, Ron. Intraempreendedorismo na prática: um guia de inovação. Campus: 2004RIES, Eric. A startup enxuta. Leya Editora. 1ª edição, 2011SANTOS. S.A. e CUNHA, N.C.V (orgs.). Empresas de Base Tecnológica: Conceitos, instrumentos e recursos. Unicorpore, 2005THIEL, Peter. De Zero a UM: O que aprender sobre empreendedorismo com Vale do Silício. Objetiva. 1ª edição, 2014TIMMONS; Jeffry; DORNELAS, José. SPINELLI, Stephen. A criação de novos negócios – empreendedorismo para o século 21. Editora Campus. 2010.", $false, $false, $false, $false, $false, $true, 1, $false, "BLANK, Steve Gary. Do Sonho a realização em 4 passos: Estratégias para a criação de empresas de sucesso. Editora Evora. 3ª edição, 2008^lBLANK, Steve; DORF, Bob. STARTUP: Manual do Empreendedorismo. O guia passo a passo para construir uma grande empresa. Alta Books Editora.  1ª edição, 2014.^lCECCONELO, Antonio; AJZENTAL, Alberto. A construção do plano de negócios. Ed. Saraiva, 1ª edição, 2008.^lCHIAVENATO, Idalberto. Empreendedorismo – dando asas ao espírito empreendedor. Ed. Saraiva, 3ª edição, 2008.^lDOLABELA, Fernando. O Segredo de Luísa. Rio de Janeiro: Sextante, 2008. ^lDORNELAS, Jose. Empreendedorismo: transformando ideias em negócios. Editora Campus. 1ª edição, 2001^lDORNELAS, Jose. Empreendedorismo na prática. LTC. 3ª edição, 2015^lDORNELAS, Jose Carlos Assis. Empreendedorismo na prática – mitos e verdades do empreendedor de sucesso. Elsevier/Campus: Rio de Janeiro, 2007. ^lFILION, L. J.; Visão e Relações: Elementos para um Metamodelo da Atividade Empreendedora. International Small Business Journal, 1991. Tradução de Costa, S.R. ^lFILION, L. J.; - O planejamento do seu Sistema de Aprendizagem Empresarial: Identifique uma Visão e Avalie o seu Sistema de Relações. Revista de Administração de Empresas, FGV, São Paulo, jul/set. 1991, pag. 31(3): 63:71. ^lHASHIMOTO, Marcos. Espírito empreendedor nas organizações – aumentando a competitividade através do intraempreendedorismo. São Paulo: Saraiva, 2006. ^lHISRICH, Robert; PETERS, Michael.  Empreendedorismo. 5.ed. - Porto Alegre: Bookman, 2004. ^lOSTERWALDER, Alexander. Inovação Em Modelos de Negócios – Business Model Generation. Editora Alta Books, 2011^lPINCHOT, Gifford; PELLMAN, Ron. Intraempreendedorismo na prática: um guia de inovação. Campus: 2004^lRIES, Eric. A startup enxuta. Leya Editora. 1ª edição, 2011^lSANTOS. S.A. e CUNHA, N.C.V (orgs.). Empresas de Base Tecnológica: Conceitos, instrumentos e recursos. Unicorpore, 2005^lTHIEL, Peter. De Zero a UM: O que aprender sobre empreendedorismo com Vale do Silício. Objetiva. 1ª edição, 2014^lTIMMONS; Jeffry; DORNELAS, José. SPINELLI, Stephen. A criação de novos negócios – empreendedorismo para o século 21. Editora Campus. 2010.", 2)
Write-Host "Target 2 found:" $found2

# Target 3
$found3 = $find.Execute("1. Characteristics of Entrepreneurial Behavior: Search for opportunities and initiative. Take calculated risks. Demand for quality and efficiency. Persistence. Commitment. Information search. Setting goals. Systematic monitoring and planning. Persuasion and networking. Independence and self-confidence.2. Strategy, Innovation and marketing.3. Design Thinking.4. Business Model (Business Model Canvas and Lean Startup - Lean Canvas): Problem. Customer Segment. Unique Value Proposition. Solution. Key Metrics. Channels. Cost Structure. Revenue Streams. Unfair Advantage.5. Minimum Viable Product: Build-Measure-Learn Cycle. Customer lifetime value. Rapid prototyping.6. Process Management and Agile Project Management7. Business Plan: Marketing, Finance, Human Resources, Product Development and Information and Communication Technology. 8. Proposal for the creation of a startup, from the business model to the business plan, including the assembly of the minimum viable product and a PITCH round. 9. Development of practical extension activity (production of digital content on entrepreneurship and innovation)10. Visit (complementary didactic trip) to an environment of innovation and entrepreneurship (eg incubator/accelerator or technology park), to understand the development of entrepreneurship and innovation processes.", $false, $false, $false, $false, $false, $true, 1, $false, "1. Characteristics of Entrepreneurial Behavior: Search for opportunities and initiative. Take calculated risks. Demand for quality and efficiency. Persistence. Commitment. Information search. Setting goals. Systematic monitoring and planning. Persuasion and networking. Independence and self-confidence.^l2. Strategy, Innovation and marketing.^l3. Design Thinking.^l4. Business Model (Business Model Canvas and Lean Startup - Lean Canvas): Problem. Customer Segment. Unique Value Proposition. Solution. Key Metrics. Channels. Cost Structure. Revenue Streams. Unfair Advantage.^l5. Minimum Viable Product: Build-Measure-Learn Cycle. Customer lifetime value. Rapid prototyping.^l6. Process Management and Agile Project Management^l7. Business Plan: Marketing, Finance, Human Resources, Product Development and Information and Communication Technology. ^l8. Proposal for the creation of a startup, from the business model to the business plan, including the assembly of the minimum viable product and a PITCH round. ^l9. Development of practical extension activity (production of digital content on entrepreneurship and innovation)^l10. Visit (complementary didactic trip) to an environment of innovation and entrepreneurship (eg incubator/accelerator or technology park), to understand the development of entrepreneurship and innovation processes.", 2)
Write-Host "Target 3 found:" $found3

# Target 4
$found4 = $find.Execute("São objetivos da atividade Extensionista:- Disseminar a cultura empreendedora e inovadora aos estudantes de Ensino Médio;- Estimular os estudantes de Ensino Médio para o desenvolvimento de sua capacidade empreendedora, a busca de oportunidades, a geração do autoemprego e o desenvolvimento de atitudes empreendedoras e criativas.", $false, $false, $false, $false, $false, $true, 1, $false, "São objetivos da atividade Extensionista:^l- Disseminar a cultura empreendedora e inovadora aos estudantes de Ensino Médio;^l- Estimular os estudantes de Ensino Médio para o desenvolvimento de sua capacidade empreendedora, a busca de oportunidades, a geração do autoemprego e o desenvolvimento de atitudes empreendedoras e criativas.", 2)
Write-Host "Target 4 found:" $found4

# Target 5
$found5 = $find.Execute("Esta atividade é denominada Engenharia e Negócios – Oficina de Empreendedorismo e Inovação.A atividade consiste na realização de uma oficina de Empreendedorismo e Inovação com estudantes do Ensino Médio. Tal oficina poderá ocorrer em escolas de ensino médio ou em organizações sociais ou representativas das comunidades da cidade de Lorena (ex: associações de bairros).Etapas:1.Planejamento da Oficina: definição dos temas (na área de empreendedorismo e inovação) a serem desenvolvidos, que pode incluir concursos de ideias, desafios de negócios, exposição de conteúdos, entre outras atividades, com a consequente preparação dos materiais (slides, vídeos, exercícios etc.) que serão utilizados nas oficinas. Os alunos serão os responsáveis por agendar a realização das oficinas com as escolas ou com outras organizações previamente aprovadas pelo professor da disciplina.2.Realização da Oficina: poderá ser aplicada em uma ou duas etapas (dias diferentes), somando no mínimo 4 horas totais de aplicação.3.Preparação de Relatos em Vídeo: criação de um vídeo relatando o desenvolvimento da oficina, com os aprendizados adquiridos, para ser disponibilizado para a comunidade.4.Autoavaliação pelo Grupo: avaliar os resultados da avaliação da atividade aplicada aos estudantes do ensino médio, para identificar o aprendizado e os pontos a melhorar para as próximas oficinas.", $false, $false, $false, $false, $false, $true, 1, $false, "Esta atividade é denominada Engenharia e Negócios – Oficina de Empreendedorismo e Inovação.^lA atividade consiste na realização de uma oficina de Empreendedorismo e Inovação com estudantes do Ensino Médio. Tal oficina poderá ocorrer em escolas de ensino médio ou em organizações sociais ou representativas das comunidades da cidade de Lorena (ex: associações de bairros).^lEtapas:^l1.Planejamento da Oficina: definição dos temas (na área de empreendedorismo e inovação) a serem desenvolvidos, que pode incluir concursos de ideias, desafios de negócios, exposição de conteúdos, entre outras atividades, com a consequente preparação dos materiais (slides, vídeos, exercícios etc.) que serão utilizados nas oficinas. Os alunos serão os responsáveis por agendar a realização das oficinas com as escolas ou com outras organizações previamente aprovadas pelo professor da disciplina.^l2.Realização da Oficina: poderá ser aplicada em uma ou duas etapas (dias diferentes), somando no mínimo 4 horas totais de aplicação.^l3.Preparação de Relatos em Vídeo: criação de um vídeo relatando o desenvolvimento da oficina, com os aprendizados adquiridos, para ser disponibilizado para a comunidade.^l4.Autoavaliação pelo Grupo: avaliar os resultados da avaliação da atividade aplicada aos estudantes do ensino médio, para identificar o aprendizado e os pontos a melhorar para as próximas oficinas.", 2)
Write-Host "Target 5 found:" $found5
